# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Updates DeliveryPlan (sheet1) ori_deployment_uid/material/delivery_qty/VFR values
# and VehicleLog (sheet2) total_units/total_volume/VFR values to reflect the
# corrected future-production-plan / available-inventory calculation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# DeliveryPlan sheet
# ---------------------------------------------------------------------------
$delivery = $wb.Worksheets.Item("DeliveryPlan")

# Row 2
$delivery.Range("B2").Value = "MAT_A|PLANT_001|DC_001|2024-01-02|net demand for forecast|000014"
$delivery.Range("C2").Value = "MAT_A"
$delivery.Range("I2").Value = 20
$delivery.Range("M2").Value = 0.9625

# Row 3
$delivery.Range("B3").Value = "MAT_B|PLANT_001|DC_001|2024-01-02|net demand for forecast|000019"
$delivery.Range("C3").Value = "MAT_B"
$delivery.Range("I3").Value = 40
$delivery.Range("M3").Value = 0.9625

# Row 4
$delivery.Range("B4").Value = "MAT_A|PLANT_001|DC_001|2024-01-03|net demand for forecast|000015"
$delivery.Range("C4").Value = "MAT_A"
$delivery.Range("I4").Value = 20
$delivery.Range("M4").Value = 0.9625

# Row 5
$delivery.Range("B5").Value = "MAT_B|PLANT_001|DC_001|2024-01-03|net demand for forecast|000020"
$delivery.Range("C5").Value = "MAT_B"
$delivery.Range("I5").Value = 35
$delivery.Range("M5").Value = 0.9625

# Row 6
$delivery.Range("B6").Value = "MAT_B|PLANT_001|DC_001|2024-01-03|net demand for forecast|000020"
$delivery.Range("C6").Value = "MAT_B"
$delivery.Range("I6").Value = 5
$delivery.Range("M6").Value = 0.9625

# Row 7
$delivery.Range("B7").Value = "MAT_A|PLANT_001|DC_001|2024-01-04|net demand for forecast|000016"
$delivery.Range("C7").Value = "MAT_A"
$delivery.Range("I7").Value = 20
$delivery.Range("M7").Value = 0.9625

# Row 8
$delivery.Range("B8").Value = "MAT_B|PLANT_001|DC_001|2024-01-04|net demand for forecast|000021"
$delivery.Range("C8").Value = "MAT_B"
$delivery.Range("I8").Value = 40
$delivery.Range("M8").Value = 0.9625

# Row 9
$delivery.Range("B9").Value = "MAT_A|PLANT_001|DC_001|2024-01-05|net demand for forecast|000017"
$delivery.Range("C9").Value = "MAT_A"
$delivery.Range("I9").Value = 20
$delivery.Range("M9").Value = 0.9625

# Row 10
$delivery.Range("B10").Value = "MAT_B|PLANT_001|DC_001|2024-01-05|net demand for forecast|000022"
$delivery.Range("C10").Value = "MAT_B"
$delivery.Range("I10").Value = 30
$delivery.Range("M10").Value = 0.9625

# Row 11 (only the ori_deployment_uid sequence number changes)
$delivery.Range("B11").Value = "MAT_A|PLANT_001|DC_002|2024-01-02|net demand for forecast|000010"

# Row 12
$delivery.Range("B12").Value = "MAT_A|PLANT_001|DC_002|2024-01-03|net demand for forecast|000011"

# Row 13
$delivery.Range("B13").Value = "MAT_A|PLANT_001|DC_002|2024-01-03|net demand for forecast|000011"

# Row 14
$delivery.Range("B14").Value = "MAT_A|PLANT_001|DC_002|2024-01-04|net demand for forecast|000012"

# ---------------------------------------------------------------------------
# VehicleLog sheet
# ---------------------------------------------------------------------------
$vehicleLog = $wb.Worksheets.Item("VehicleLog")

# Row 2
$vehicleLog.Range("G2").Value = 115
$vehicleLog.Range("I2").Value = 192.5
$vehicleLog.Range("K2").Value = 0.9625

# Row 3
$vehicleLog.Range("G3").Value = 115
$vehicleLog.Range("I3").Value = 192.5
$vehicleLog.Range("K3").Value = 0.9625
